# Rename the worksheet tab from "Session" to "Anatomy"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Anatomy"

# New QR-scanner log rows to append (rows 38-76)
$data = @(
    @("234612","Anatomy","17/12/2025","10:21:53","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234302","Anatomy","17/12/2025","10:21:56","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234751","Anatomy","17/12/2025","10:21:59","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234314","Anatomy","17/12/2025","10:22:01","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234609","Anatomy","17/12/2025","10:22:04","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234319","Anatomy","17/12/2025","10:22:06","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234281","Anatomy","17/12/2025","10:22:09","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234271","Anatomy","17/12/2025","10:22:11","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234316","Anatomy","17/12/2025","10:22:14","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234313","Anatomy","17/12/2025","10:22:17","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234330","Anatomy","17/12/2025","10:22:20","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234361","Anatomy","17/12/2025","10:22:23","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234365","Anatomy","17/12/2025","10:22:25","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234362","Anatomy","17/12/2025","10:22:28","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234194","Anatomy","17/12/2025","10:22:33","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234221","Anatomy","17/12/2025","10:22:35","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234260","Anatomy","17/12/2025","10:22:38","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234480","Anatomy","17/12/2025","10:22:39","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234317","Anatomy","17/12/2025","10:22:40","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234259","Anatomy","17/12/2025","10:22:47","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234284","Anatomy","17/12/2025","10:22:50","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234198","Anatomy","17/12/2025","10:38:34","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234199","Anatomy","17/12/2025","10:38:37","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234441","Anatomy","17/12/2025","10:53:31","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234176","Anatomy","17/12/2025","10:53:35","Scan","nahla.nagiub@med.asu.edu.eg"),
    @("234607","Anatomy","17/12/2025","10:53:54","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234294","Anatomy","17/12/2025","10:54:02","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234878","Anatomy","17/12/2025","10:54:12","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234246","Anatomy","17/12/2025","10:54:32","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234293","Anatomy","17/12/2025","10:54:41","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234382","Anatomy","17/12/2025","10:55:06","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234389","Anatomy","17/12/2025","10:55:21","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234522","Anatomy","17/12/2025","10:55:28","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234508","Anatomy","17/12/2025","10:55:36","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234523","Anatomy","17/12/2025","10:55:44","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234277","Anatomy","17/12/2025","10:56:10","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("231068","Anatomy","17/12/2025","10:56:18","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234286","Anatomy","17/12/2025","10:56:30","Manual","nahla.nagiub@med.asu.edu.eg"),
    @("234266","Anatomy","17/12/2025","10:56:42","Manual","nahla.nagiub@med.asu.edu.eg")
)

$startRow = 38
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    # Column A (Student ID) looks numeric - prefix with an apostrophe so Excel
    # keeps it stored as text, matching the rest of the log sheet.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

Write-Output "Appended $($data.Count) rows; sheet renamed to $($ws.Name)"
